$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: four new fixture headers (AV1:AY1) ---
# Clone AU1's header style (bold, centered, thin border) onto the new header cells
# before writing their text, so they reuse the existing style index (s="1") exactly
# like the rest of row 1, instead of Excel allocating a brand new style record.
$ws.Range("AU1").Copy()
$ws.Range("AV1:AY1").PasteSpecial(-4122)

$ws.Range("AV1").Value = 'Spain vs Germany'
$ws.Range("AW1").Value = 'England vs Switzerland'
$ws.Range("AX1").Value = 'Netherlands vs Turkey'
$ws.Range("AY1").Value = 'Portugal vs France'

# --- Rows 2-33: prediction cell for each of the 4 new fixtures ---
# Most users left these fixtures blank; the source workbook still records an explicit
# (empty) inline-string cell for every user row rather than leaving the cell absent, so
# a lone leading apostrophe is used to force an empty *text* cell (Excel's literal-text
# quote prefix) instead of Value = "" which the engine treats as "no cell at all".

# Row 2
$ws.Range("AV2").Value = "'"
$ws.Range("AW2").Value = "'"
$ws.Range("AX2").Value = "'"
$ws.Range("AY2").Value = "'"
# Row 3
$ws.Range("AV3").Value = "'"
$ws.Range("AW3").Value = "'"
$ws.Range("AX3").Value = "'"
$ws.Range("AY3").Value = "'"
# Row 4
$ws.Range("AV4").Value = "'"
$ws.Range("AW4").Value = "'"
$ws.Range("AX4").Value = "'"
$ws.Range("AY4").Value = "'"
# Row 5
$ws.Range("AV5").Value = "'"
$ws.Range("AW5").Value = "'"
$ws.Range("AX5").Value = "'"
$ws.Range("AY5").Value = "'"
# Row 6
$ws.Range("AV6").Value = "'"
$ws.Range("AW6").Value = '[''England'', 2, 1]'
$ws.Range("AX6").Value = '[''Draw'', 1, 1]'
$ws.Range("AY6").Value = "'"
# Row 7
$ws.Range("AV7").Value = "'"
$ws.Range("AW7").Value = '[''England'', 2, 0]'
$ws.Range("AX7").Value = '[''Netherlands'', 2, 0]'
$ws.Range("AY7").Value = '[''France'', 0, 2]'
# Row 8
$ws.Range("AV8").Value = "'"
$ws.Range("AW8").Value = "'"
$ws.Range("AX8").Value = "'"
$ws.Range("AY8").Value = "'"
# Row 9
$ws.Range("AV9").Value = '[''Draw'', 1, 1]'
$ws.Range("AW9").Value = '[''Draw'', 1, 1]'
$ws.Range("AX9").Value = '[''Netherlands'', 2, 1]'
$ws.Range("AY9").Value = '[''Draw'', 1, 1]'
# Row 10
$ws.Range("AV10").Value = '[''Spain'', 2, 1]'
$ws.Range("AW10").Value = '[''Draw'', 0, 0]'
$ws.Range("AX10").Value = '[''Netherlands'', 2, 0]'
$ws.Range("AY10").Value = '[''France'', 0, 1]'
# Row 11
$ws.Range("AV11").Value = "'"
$ws.Range("AW11").Value = "'"
$ws.Range("AX11").Value = "'"
$ws.Range("AY11").Value = "'"
# Row 12
$ws.Range("AV12").Value = "'"
$ws.Range("AW12").Value = "'"
$ws.Range("AX12").Value = "'"
$ws.Range("AY12").Value = "'"
# Row 13
$ws.Range("AV13").Value = "'"
$ws.Range("AW13").Value = "'"
$ws.Range("AX13").Value = "'"
$ws.Range("AY13").Value = "'"
# Row 14
$ws.Range("AV14").Value = '[''Germany'', 0, 1]'
$ws.Range("AW14").Value = '[''England'', 2, 1]'
$ws.Range("AX14").Value = '[''Turkey'', 0, 1]'
$ws.Range("AY14").Value = '[''Portugal'', 2, 1]'
# Row 15
$ws.Range("AV15").Value = "'"
$ws.Range("AW15").Value = "'"
$ws.Range("AX15").Value = "'"
$ws.Range("AY15").Value = "'"
# Row 16
$ws.Range("AV16").Value = "'"
$ws.Range("AW16").Value = "'"
$ws.Range("AX16").Value = "'"
$ws.Range("AY16").Value = "'"
# Row 17
$ws.Range("AV17").Value = "'"
$ws.Range("AW17").Value = "'"
$ws.Range("AX17").Value = "'"
$ws.Range("AY17").Value = "'"
# Row 18
$ws.Range("AV18").Value = "'"
$ws.Range("AW18").Value = "'"
$ws.Range("AX18").Value = "'"
$ws.Range("AY18").Value = "'"
# Row 19
$ws.Range("AV19").Value = '[''Draw'', 2, 2]'
$ws.Range("AW19").Value = '[''England'', 2, 1]'
$ws.Range("AX19").Value = '[''Netherlands'', 2, 0]'
$ws.Range("AY19").Value = '[''France'', 0, 2]'
# Row 20
$ws.Range("AV20").Value = "'"
$ws.Range("AW20").Value = "'"
$ws.Range("AX20").Value = "'"
$ws.Range("AY20").Value = "'"
# Row 21
$ws.Range("AV21").Value = "'"
$ws.Range("AW21").Value = "'"
$ws.Range("AX21").Value = "'"
$ws.Range("AY21").Value = "'"
# Row 22
$ws.Range("AV22").Value = "'"
$ws.Range("AW22").Value = "'"
$ws.Range("AX22").Value = "'"
$ws.Range("AY22").Value = "'"
# Row 23
$ws.Range("AV23").Value = '[''Draw'', 1, 1]'
$ws.Range("AW23").Value = '[''England'', 1, 0]'
$ws.Range("AX23").Value = '[''Netherlands'', 2, 0]'
$ws.Range("AY23").Value = '[''France'', 0, 1]'
# Row 24
$ws.Range("AV24").Value = "'"
$ws.Range("AW24").Value = "'"
$ws.Range("AX24").Value = "'"
$ws.Range("AY24").Value = "'"
# Row 25
$ws.Range("AV25").Value = "'"
$ws.Range("AW25").Value = "'"
$ws.Range("AX25").Value = "'"
$ws.Range("AY25").Value = "'"
# Row 26
$ws.Range("AV26").Value = "'"
$ws.Range("AW26").Value = "'"
$ws.Range("AX26").Value = "'"
$ws.Range("AY26").Value = "'"
# Row 27
$ws.Range("AV27").Value = "'"
$ws.Range("AW27").Value = "'"
$ws.Range("AX27").Value = "'"
$ws.Range("AY27").Value = "'"
# Row 28
$ws.Range("AV28").Value = "'"
$ws.Range("AW28").Value = "'"
$ws.Range("AX28").Value = "'"
$ws.Range("AY28").Value = "'"
# Row 29
$ws.Range("AV29").Value = "'"
$ws.Range("AW29").Value = "'"
$ws.Range("AX29").Value = "'"
$ws.Range("AY29").Value = "'"
# Row 30
$ws.Range("AV30").Value = "'"
$ws.Range("AW30").Value = "'"
$ws.Range("AX30").Value = "'"
$ws.Range("AY30").Value = "'"
# Row 31
$ws.Range("AV31").Value = '[''Spain'', 2, 1]'
$ws.Range("AW31").Value = '[''Switzerland'', 1, 2]'
$ws.Range("AX31").Value = '[''Netherlands'', 3, 1]'
$ws.Range("AY31").Value = '[''France'', 0, 1]'
# Row 32
$ws.Range("AV32").Value = '[''Germany'', 2, 3]'
$ws.Range("AW32").Value = '[''England'', 2, 0]'
$ws.Range("AX32").Value = '[''Netherlands'', 3, 0]'
$ws.Range("AY32").Value = '[''France'', 0, 3]'
# Row 33
$ws.Range("AV33").Value = "'"
$ws.Range("AW33").Value = "'"
$ws.Range("AX33").Value = "'"
$ws.Range("AY33").Value = "'"
